# Add a new worksheet "Summary_Epsilon_5000" at the end of the workbook
# containing a two-column (Attribute, Count) summary table.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Summary_Epsilon_5000"

# Header row
$headers = @("Attribute", "Count")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Style the header row to match the other sheets: bold, thin box border,
# centered horizontally, top-aligned vertically.
$headerRange = $ws.Range("A1:B1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Data rows: Attribute, Count
$data = @(
    @("HDI", 1925),
    @("GDP", 1673),
    @("GINI", 1626),
    @("Continent", 1615),
    @("Student", 1591),
    @("SexualOrientation", 1535),
    @("Hobby", 1479),
    @("Gender", 1339),
    @("RaceEthnicity", 1099),
    @("Country", 1056),
    @("UndergradMajor", 1020),
    @("Dependents", 971),
    @("FormalEducation", 833),
    @("DevType", 608),
    @("Age", 565),
    @("HoursComputer", 475),
    @("Exercise", 154),
    @("EducationParents", 49),
    @("YearsCoding", 18)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}

$ws.Range("A1").Select() | Out-Null
